$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple D-column price/volume text updates (stored as text, not numbers)
$updates = @{
    2  = "271.61"
    4  = "6.342"
    5  = "0.06209"
    6  = "3.652"
    9  = "0.8304"
    10 = "0.01380"
    11 = "0.1603"
    12 = "0.08295"
    13 = "0.03431"
    14 = "0.03174"
    15 = "0.09337"
    16 = "3.860"
    17 = "0.001632"
    18 = "0.04731"
    19 = "0.006321"
    20 = "0.005671"
    21 = "0.001077"
    22 = "0.0001500"
    23 = "3.720"
    24 = "2.325"
    25 = "0.3347"
    27 = "0.0002704"
    40 = "0.04690"
    41 = "0.007032"
    44 = "0.01166"
    45 = "0.00006245"
    46 = "0.0009908"
    48 = "0.9202"
    49 = "0.002087"
    50 = "0.00001400"
    51 = "0.01240"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
}

# Row 42/43: the two coins swap places, with new price values
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003801"
$ws.Range("E42").Value = "41CEJICEJIWorstin24h"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1162"
$ws.Range("E43").Value = "42BKEXTokenBKK"
